{"js": "// Office.js (Word JavaScript API) script\n// Applies the resume restructuring edit described by the diff:\n//  - Simplifies header formatting (name size, removes colors, merges contact lines)\n//  - Removes the \"Professional Title\" placeholder paragraph\n//  - Strips heading color overrides\n//  - Merges competency category headers into their body paragraphs\n//  - Expands the single placeholder job entry into the full 9-job work history\n//    (using Heading3 for job titles instead of bold/colored runs, plain bullet\n//    paragraphs with \"\u2022 \" instead of ListBullet-styled \"\u25b8 \"/\"\u2713 \" paragraphs)\n//  - Converts \"Geospatial Platform Engineering\" to Heading3 and its bullets to\n//    plain \"\u2022 \" paragraphs\n//  - Updates the page margins (top/bottom 864->1440, left/right 864->1800 twips)\n\nconst body = context.document.body;\n\nconst newBodyXml = `    <w:p>\n      <w:pPr>\n        <w:jc w:val=\"center\"/>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:b/>\n          <w:sz w:val=\"28\"/>\n        </w:rPr>\n        <w:t>Dheeraj Chand</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:jc w:val=\"center\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>PROFESSIONAL SUMMARY</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Senior Data Engineer with 21 years of expertise in geospatial data platforms, big data processing, and distributed systems architecture. Deep specialist in Apache Spark/Sedona for large-scale geospatial analytics, with fluency across ESRI, OSGeo, and SAFE FME technology stacks. Proven track record architecting production systems serving thousands of users, implementing PySpark pipelines processing billions of spatial records, and leading engineering teams. Expert in full-stack geospatial development from PostGIS database optimization to React-based mapping interfaces.</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>CORE COMPETENCIES</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Big Data &amp; Geospatial Processing: Apache Spark: PySpark, Spark SQL, Scala Spark, Sedona (geospatial), distributed processing \u2022 Geospatial Databases: PostGIS (advanced), Oracle Spatial, spatial indexing, query optimization \u2022 ETL/ELT: dbt, Informatica, CDAP, custom PySpark pipelines, data governance frameworks \u2022 Cloud Platforms: AWS (EC2, RDS, S3), Snowflake, Hadoop clusters, distributed computing \u2022 Streaming: Real-time data processing, Kafka integration, event-driven architectures</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>GIS Technology Stack: ESRI: ArcGIS Server, ArcGIS Pro, enterprise geodatabases, ModelBuilder, ArcPy scripting \u2022 OSGeo: QGIS, GRASS GIS, GDAL/OGR, GeoServer, spatial analysis workflows \u2022 SAFE FME: Data transformation, format conversion, spatial ETL, enterprise integration \u2022 Web Mapping: OpenLayers, Leaflet, MapBox, tile servers, WMS/WFS services \u2022 Spatial Analysis: Clustering algorithms, boundary estimation, network analysis, geostatistics</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Software Development &amp; Architecture: Python: Django/GeoDjango, Flask, Pandas, NumPy, SciKit-Learn, spatial libraries \u2022 JVM: Scala (Spark), Java (GeoTools, enterprise), Groovy scripting \u2022 Web Technologies: React, JavaScript, d3.js, RESTful APIs, microservices \u2022 Databases: PostgreSQL/PostGIS, Oracle, MySQL, MongoDB, spatial optimization \u2022 DevOps: Docker, Kubernetes, CI/CD (GitLab, GitHub), Airflow, Celery, nginx</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>PROFESSIONAL EXPERIENCE</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>PARTNER &amp; SENIOR DATA ARCHITECT - Siege Analytics, Washington, DC | January 2014 \u2013 Present</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Geospatial Data Platform Architecture and Big Data Engineering</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Architected and engineered production geospatial platforms serving thousands of analysts</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built enterprise-scale ETL pipelines using PySpark and Sedona processing billions of geospatial records with advanced spatial clustering algorithms</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed multi-tenant data warehouse integrating Census, electoral, and demographic data using PostGIS and Spark SQL optimization</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Implemented fraud detection systems processing multi-terabyte datasets with real-time spatial analysis capabilities</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Created parametric boundary estimation algorithms using PostGIS and GRASS without machine learning dependencies</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Led integration of ESRI ArcGIS Server, OSGeo tools (QGIS, GRASS), and SAFE FME for enterprise geospatial workflows</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>PRINCIPAL TECHNICAL ARCHITECT - Clarity and Rigour, Washington, DC | 2012 \u2013 2014</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Geospatial Systems Architecture and Development</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Architected geospatial analysis frameworks and mapping applications for electoral research</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed custom visualization tools and interactive dashboards using JavaScript and OpenLayers</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Designed data processing pipelines for large-scale demographic and geographic datasets</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Implemented PostGIS spatial databases and optimized geospatial query performance</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>DIRECTOR OF TECHNOLOGY - Helm, Washington, DC | 2010 \u2013 2012</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Data Platform Architecture and Engineering Leadership</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Led technical architecture and development of data-driven political technology platforms</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Designed and implemented scalable data platforms using Python, Django, and PostgreSQL</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built RESTful APIs and microservices architecture for campaign data integration</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Managed engineering teams and established development best practices and CI/CD pipelines</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>SENIOR TECHNICAL ANALYST - GSD&amp;M, Austin, TX | 2008 \u2013 2010</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Campaign Technology and Data Engineering</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed voter targeting models and demographic analysis tools using Python and R</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built web applications and data visualization systems for campaign analytics</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Created data integration systems connecting multiple campaign data sources</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Implemented machine learning algorithms and statistical models for voter behavior prediction</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>TECHNICAL COORDINATOR - Progressive Change Campaign Committee, Washington, DC | 2006 \u2013 2008</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Political Technology and Data Systems</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Coordinated technical operations and data systems for political campaigns</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed data collection and analysis protocols for campaign research</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built custom applications and tools for voter engagement and campaign management</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Supported technical infrastructure and data processing for progressive political initiatives</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>SOFTWARE ENGINEER - Salsa Labs, Inc., Washington, DC | 2004 \u2013 2006</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Political Technology Development</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed software solutions for political campaigns using PHP, JavaScript, and MySQL</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built web applications for voter engagement and campaign management</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Integrated third-party APIs and data sources for campaign tools</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Collaborated with political strategists to translate requirements into technical solutions</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 \u2013 2004</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Nonprofit Technology Integration and Development</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed data management systems and web applications for social justice organizations</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built custom applications for community engagement using PHP, MySQL, and web technologies</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Provided technical training and support to nonprofit staff</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Integrated technology solutions within organizational frameworks for advocacy work</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>PROGRAMMER - Lake Research Partners, Washington, DC | 2001 \u2013 2002</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Political Research and Data Analysis Tools</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed data analysis tools for political polling and research using Python and R</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built statistical models and data visualization tools for research presentations</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Created automated reporting systems and data processing pipelines for survey analysis</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Supported senior researchers with technical analysis and data processing automation</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>TECHNICAL COORDINATOR - The Feldman Group, Washington, DC | 2000 \u2013 2001</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Field Operations Technology and Data Management</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed data collection and management systems for political field operations</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built databases and reporting tools for campaign field work and voter outreach</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Trained field staff on data collection protocols and quality control systems</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Analyzed field data using statistical methods to inform campaign strategy and research findings</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>KEY ACHIEVEMENTS AND IMPACT</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>Geospatial Platform Engineering</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Architected redistricting platform processing Census data for thousands of analysts with real-time PostGIS collaborative editing</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built boundary estimation system using advanced PostGIS algorithms and incomplete data without machine learning requirements</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed geospatial simulation platform integrating multi-agent modeling with web interface</w:t>\n      </w:r>\n    </w:p>\n`;\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n${newBodyXml}</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\n// Replace the entire body content with the restructured paragraphs while\n// keeping the existing sectPr (page size, header/footer distances, etc.)\nbody.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// Update the page margins via the document-level PageSetup object\n// (864 twips = 43.2pt -> 1440 twips = 72pt; 864 twips -> 1800 twips = 90pt)\nconst pageSetup = context.document.pageSetup;\npageSetup.topMargin = 72;\npageSetup.bottomMargin = 72;\npageSetup.leftMargin = 90;\npageSetup.rightMargin = 90;\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Applies the resume restructuring edit described by the diff:\n#  - Simplifies header formatting (name size, removes colors, merges contact lines)\n#  - Removes the \"Professional Title\" placeholder paragraph\n#  - Strips heading color overrides\n#  - Merges competency category headers into their body paragraphs\n#  - Expands the single placeholder job entry into the full 9-job work history\n#    (using Heading3 for job titles instead of bold/colored runs, plain bullet\n#    paragraphs with bullet-character prefixes instead of ListBullet-styled ones)\n#  - Converts \"Geospatial Platform Engineering\" to Heading3 and its bullets to\n#    plain bullet-character paragraphs\n#  - Updates the page margins (top/bottom 864->1440, left/right 864->1800 twips)\n\n$d = $word.ActiveDocument\n\n$newBodyXml = @'\n    <w:p>\n      <w:pPr>\n        <w:jc w:val=\"center\"/>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:b/>\n          <w:sz w:val=\"28\"/>\n        </w:rPr>\n        <w:t>Dheeraj Chand</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:jc w:val=\"center\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>PROFESSIONAL SUMMARY</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Senior Data Engineer with 21 years of expertise in geospatial data platforms, big data processing, and distributed systems architecture. Deep specialist in Apache Spark/Sedona for large-scale geospatial analytics, with fluency across ESRI, OSGeo, and SAFE FME technology stacks. Proven track record architecting production systems serving thousands of users, implementing PySpark pipelines processing billions of spatial records, and leading engineering teams. Expert in full-stack geospatial development from PostGIS database optimization to React-based mapping interfaces.</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>CORE COMPETENCIES</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Big Data &amp; Geospatial Processing: Apache Spark: PySpark, Spark SQL, Scala Spark, Sedona (geospatial), distributed processing \u2022 Geospatial Databases: PostGIS (advanced), Oracle Spatial, spatial indexing, query optimization \u2022 ETL/ELT: dbt, Informatica, CDAP, custom PySpark pipelines, data governance frameworks \u2022 Cloud Platforms: AWS (EC2, RDS, S3), Snowflake, Hadoop clusters, distributed computing \u2022 Streaming: Real-time data processing, Kafka integration, event-driven architectures</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>GIS Technology Stack: ESRI: ArcGIS Server, ArcGIS Pro, enterprise geodatabases, ModelBuilder, ArcPy scripting \u2022 OSGeo: QGIS, GRASS GIS, GDAL/OGR, GeoServer, spatial analysis workflows \u2022 SAFE FME: Data transformation, format conversion, spatial ETL, enterprise integration \u2022 Web Mapping: OpenLayers, Leaflet, MapBox, tile servers, WMS/WFS services \u2022 Spatial Analysis: Clustering algorithms, boundary estimation, network analysis, geostatistics</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Software Development &amp; Architecture: Python: Django/GeoDjango, Flask, Pandas, NumPy, SciKit-Learn, spatial libraries \u2022 JVM: Scala (Spark), Java (GeoTools, enterprise), Groovy scripting \u2022 Web Technologies: React, JavaScript, d3.js, RESTful APIs, microservices \u2022 Databases: PostgreSQL/PostGIS, Oracle, MySQL, MongoDB, spatial optimization \u2022 DevOps: Docker, Kubernetes, CI/CD (GitLab, GitHub), Airflow, Celery, nginx</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>PROFESSIONAL EXPERIENCE</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>PARTNER &amp; SENIOR DATA ARCHITECT - Siege Analytics, Washington, DC | January 2014 \u2013 Present</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Geospatial Data Platform Architecture and Big Data Engineering</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Architected and engineered production geospatial platforms serving thousands of analysts</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built enterprise-scale ETL pipelines using PySpark and Sedona processing billions of geospatial records with advanced spatial clustering algorithms</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed multi-tenant data warehouse integrating Census, electoral, and demographic data using PostGIS and Spark SQL optimization</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Implemented fraud detection systems processing multi-terabyte datasets with real-time spatial analysis capabilities</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Created parametric boundary estimation algorithms using PostGIS and GRASS without machine learning dependencies</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Led integration of ESRI ArcGIS Server, OSGeo tools (QGIS, GRASS), and SAFE FME for enterprise geospatial workflows</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>PRINCIPAL TECHNICAL ARCHITECT - Clarity and Rigour, Washington, DC | 2012 \u2013 2014</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Geospatial Systems Architecture and Development</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Architected geospatial analysis frameworks and mapping applications for electoral research</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed custom visualization tools and interactive dashboards using JavaScript and OpenLayers</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Designed data processing pipelines for large-scale demographic and geographic datasets</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Implemented PostGIS spatial databases and optimized geospatial query performance</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>DIRECTOR OF TECHNOLOGY - Helm, Washington, DC | 2010 \u2013 2012</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Data Platform Architecture and Engineering Leadership</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Led technical architecture and development of data-driven political technology platforms</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Designed and implemented scalable data platforms using Python, Django, and PostgreSQL</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built RESTful APIs and microservices architecture for campaign data integration</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Managed engineering teams and established development best practices and CI/CD pipelines</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>SENIOR TECHNICAL ANALYST - GSD&amp;M, Austin, TX | 2008 \u2013 2010</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Campaign Technology and Data Engineering</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed voter targeting models and demographic analysis tools using Python and R</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built web applications and data visualization systems for campaign analytics</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Created data integration systems connecting multiple campaign data sources</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Implemented machine learning algorithms and statistical models for voter behavior prediction</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>TECHNICAL COORDINATOR - Progressive Change Campaign Committee, Washington, DC | 2006 \u2013 2008</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Political Technology and Data Systems</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Coordinated technical operations and data systems for political campaigns</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed data collection and analysis protocols for campaign research</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built custom applications and tools for voter engagement and campaign management</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Supported technical infrastructure and data processing for progressive political initiatives</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>SOFTWARE ENGINEER - Salsa Labs, Inc., Washington, DC | 2004 \u2013 2006</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Political Technology Development</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed software solutions for political campaigns using PHP, JavaScript, and MySQL</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built web applications for voter engagement and campaign management</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Integrated third-party APIs and data sources for campaign tools</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Collaborated with political strategists to translate requirements into technical solutions</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 \u2013 2004</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Nonprofit Technology Integration and Development</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed data management systems and web applications for social justice organizations</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built custom applications for community engagement using PHP, MySQL, and web technologies</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Provided technical training and support to nonprofit staff</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Integrated technology solutions within organizational frameworks for advocacy work</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>PROGRAMMER - Lake Research Partners, Washington, DC | 2001 \u2013 2002</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Political Research and Data Analysis Tools</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed data analysis tools for political polling and research using Python and R</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built statistical models and data visualization tools for research presentations</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Created automated reporting systems and data processing pipelines for survey analysis</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Supported senior researchers with technical analysis and data processing automation</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>TECHNICAL COORDINATOR - The Feldman Group, Washington, DC | 2000 \u2013 2001</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>Field Operations Technology and Data Management</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed data collection and management systems for political field operations</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built databases and reporting tools for campaign field work and voter outreach</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Trained field staff on data collection protocols and quality control systems</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Analyzed field data using statistical methods to inform campaign strategy and research findings</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading2\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>KEY ACHIEVEMENTS AND IMPACT</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"Heading3\"/>\n      </w:pPr>\n      <w:r>\n        <w:t>Geospatial Platform Engineering</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Architected redistricting platform processing Census data for thousands of analysts with real-time PostGIS collaborative editing</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Built boundary estimation system using advanced PostGIS algorithms and incomplete data without machine learning requirements</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t>\u2022 Developed geospatial simulation platform integrating multi-agent modeling with web interface</w:t>\n      </w:r>\n    </w:p>\n\n'@\n\n$ooxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n$newBodyXml</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n\n# Replace the entire body content with the restructured paragraphs while\n# keeping the existing sectPr (page size, header/footer distances, etc.)\n$d.Content.InsertXML($ooxml)\n\n# Update the page margins via the document PageSetup object\n# (864 twips = 43.2pt -> 1440 twips = 72pt; 864 twips -> 1800 twips = 90pt)\n$ps = $d.PageSetup\n$ps.TopMargin = 72\n$ps.BottomMargin = 72\n$ps.LeftMargin = 90\n$ps.RightMargin = 90\n"}
